$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2506

$ws.Range("H80").Value = 1607.9166
$ws.Range("J80").Value = 1827.8334
$ws.Range("L80").Value = 5483.5002
$ws.Range("N80").Value = -7479.5002

$ws.Range("H83").Value = 1607.9166
$ws.Range("J83").Value = 1827.8334
$ws.Range("L83").Value = 16450.5006
$ws.Range("N83").Value = -26434.5006

$ws.Range("H113").Value = 5149.0435
$ws.Range("I113").Value = 3447.5
$ws.Range("J113").Value = 6457.923
$ws.Range("K113").Value = 3447.5
$ws.Range("L113").Value = 6457.923
$ws.Range("M113").Value = -193.5
$ws.Range("N113").Value = -12965.923

$ws.Range("H140").Value = 119938.64
$ws.Range("J140").Value = 119489.305
$ws.Range("L140").Value = 119489.305
$ws.Range("N140").Value = -129849.305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 227028.05
$ws.Range("I32").Value = 235005.77
$ws.Range("J32").Value = 55507
$ws.Range("K32").Value = 235005.77
$ws.Range("L32").Value = 55507
$ws.Range("M32").Value = -234718.77
$ws.Range("N32").Value = -56081

$ws.Range("H45").Value = 94288
$ws.Range("I45").Value = 128547.125
$ws.Range("K45").Value = 128547.125
$ws.Range("M45").Value = -128170.125

$ws.Range("H132").Value = 2591.4243
$ws.Range("J132").Value = 5976.778
$ws.Range("L132").Value = 17930.334
$ws.Range("N132").Value = -22990.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 11761
$ws.Range("I29").Value = 17035
$ws.Range("J29").Value = 3850
$ws.Range("K29").Value = 17035
$ws.Range("L29").Value = 3850
$ws.Range("M29").Value = -16746
$ws.Range("N29").Value = -4428

$ws.Range("H39").Value = 8500
$ws.Range("J39").Value = 8500
$ws.Range("L39").Value = 8500
$ws.Range("N39").Value = -9278

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3568.9546
$ws.Range("I16").Value = 1991.4166
$ws.Range("K16").Value = 1991.4166
$ws.Range("M16").Value = -1704.4166

$ws.Range("H93").Value = 27281.4
$ws.Range("I93").Value = 30703.5
$ws.Range("K93").Value = 30703.5
$ws.Range("M93").Value = -28831.5

$ws.Range("H105").Value = 2294.1428
$ws.Range("I105").Value = 1200
$ws.Range("K105").Value = 1200
$ws.Range("M105").Value = 547

$ws.Range("H107").Value = 1645.4166
$ws.Range("I107").Value = 1406.0625
$ws.Range("J107").Value = 2124.125
$ws.Range("K107").Value = 1406.0625
$ws.Range("L107").Value = 2124.125
$ws.Range("M107").Value = 513.9375
$ws.Range("N107").Value = -5964.125

$ws.Range("H113").Value = 3568.9546
$ws.Range("I113").Value = 1991.4166
$ws.Range("K113").Value = 1991.4166
$ws.Range("M113").Value = 178.5834

$ws.Range("H122").Value = 2556.875
$ws.Range("I122").Value = 2732.158
$ws.Range("K122").Value = 8196.474
$ws.Range("M122").Value = -5746.474

$ws.Range("H132").Value = 2974.2
$ws.Range("I132").Value = 2592.5715
$ws.Range("K132").Value = 7777.7145
$ws.Range("M132").Value = -5247.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 112121.22
$ws.Range("I6").Value = 168161.5
$ws.Range("J6").Value = 40.666668
$ws.Range("K6").Value = 504484.5
$ws.Range("L6").Value = 122.000004
$ws.Range("M6").Value = -504371.5
$ws.Range("N6").Value = -348.000004

$ws.Range("H16").Value = 270
$ws.Range("I16").Value = 270
$ws.Range("K16").Value = 810
$ws.Range("M16").Value = -637

$ws.Range("H25").Value = 1551.1111
$ws.Range("J25").Value = 1611
$ws.Range("L25").Value = 4833
$ws.Range("N25").Value = -5171

$ws.Range("H30").Value = 1551.1111
$ws.Range("J30").Value = 1611
$ws.Range("L30").Value = 4833
$ws.Range("N30").Value = -5037

$ws.Range("H34").Value = 3004.1667
$ws.Range("I34").Value = 208
$ws.Range("J34").Value = 5001.4287
$ws.Range("K34").Value = 624
$ws.Range("L34").Value = 15004.2861
$ws.Range("M34").Value = -540
$ws.Range("N34").Value = -15172.2861

$ws.Range("H55").Value = 4658.9375
$ws.Range("J55").Value = 5384.231
$ws.Range("L55").Value = 16152.693
$ws.Range("N55").Value = -16506.693

$ws.Range("H107").Value = 45454932
$ws.Range("J107").Value = 58823812
$ws.Range("L107").Value = 176471436
$ws.Range("N107").Value = -176475276

$ws.Range("H109").Value = 9532.666999999999
$ws.Range("J109").Value = 9532.666999999999
$ws.Range("L109").Value = 28598.001
$ws.Range("N109").Value = -30678.001

$ws.Range("H122").Value = 2073.6667
$ws.Range("I122").Value = 1110.5
$ws.Range("K122").Value = 9994.5
$ws.Range("M122").Value = -7544.5

$ws.Range("H131").Value = 4841422
$ws.Range("J131").Value = 75633.86
$ws.Range("L131").Value = 226901.58
$ws.Range("N131").Value = -236981.58

$ws.Range("H137").Value = 3862.6897
$ws.Range("J137").Value = 5954.0835
$ws.Range("L137").Value = 17862.2505
$ws.Range("N137").Value = -28062.2505

$ws.Range("H138").Value = 3831.8076
$ws.Range("I138").Value = 3147
$ws.Range("K138").Value = 9441
$ws.Range("M138").Value = -4301

$ws.Range("H139").Value = 5955640.5
$ws.Range("I139").Value = 10418522
$ws.Range("K139").Value = 31255566
$ws.Range("M139").Value = -31250426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 11481.546
$ws.Range("I41").Value = 11481.546
$ws.Range("K41").Value = 11481.546
$ws.Range("M41").Value = -11126.546

$ws.Range("H113").Value = 2451
$ws.Range("I113").Value = 2363.9167
$ws.Range("K113").Value = 2363.9167
$ws.Range("M113").Value = -193.9167000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5441.5454
$ws.Range("J46").Value = 998.5
$ws.Range("L46").Value = 998.5
$ws.Range("N46").Value = -1374.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 149999
$ws.Range("J40").Value = 149999
$ws.Range("L40").Value = 149999
$ws.Range("N40").Value = -150297

$ws.Range("H132").Value = 24550.35
$ws.Range("J132").Value = 1984.3636
$ws.Range("L132").Value = 5953.0908
$ws.Range("N132").Value = -11013.0908

$ws.Range("H136").Value = 32172.727
$ws.Range("I136").Value = 59954.234
$ws.Range("J136").Value = 2654.875
$ws.Range("K136").Value = 179862.702
$ws.Range("L136").Value = 7964.625
$ws.Range("M136").Value = -177312.702
$ws.Range("N136").Value = -13064.625

$ws.Range("H140").Value = 85734.5
$ws.Range("J140").Value = 85734.5
$ws.Range("L140").Value = 85734.5
$ws.Range("N140").Value = -96094.5
